$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 3
Write-Host ("ScrollRow after: " + $excel.ActiveWindow.ScrollRow())
Write-Host ("ScrollColumn after: " + $excel.ActiveWindow.ScrollColumn())
